$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: LP1912  (columns: A=meta, B=Hora_Scrap, C=Hora_Llegada, D=Linea,
#                             E=Minutos, F=Parada, G=Fecha)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 31/12/2025 19:21:43"
$ws1.Range("A3").Value = "Total filas: 1272"

$rows1 = @(
    @("", "19:21:32", "19:29", "16_SANTA ANA",   8,  "LP1912", "31/12/2025"),
    @("", "19:21:32", "19:29", "225_GOMEZ",      8,  "LP1912", "31/12/2025"),
    @("", "19:21:32", "19:31", "215_EL PELIGRO", 10, "LP1912", "31/12/2025"),
    @("", "19:21:32", "19:33", "23_HERNANDEZ",   12, "LP1912", "31/12/2025"),
    @("", "19:21:32", "19:41", "16_SANTA ANA",   20, "LP1912", "31/12/2025"),
    @("", "19:21:32", "19:51", "81_EL PELIGRO",  30, "LP1912", "31/12/2025"),
    @("", "19:21:32", "19:58", "14X44_ABASTO",   37, "LP1912", "31/12/2025"),
    @("", "19:21:32", "20:01", "215C_EL PATO",   40, "LP1912", "31/12/2025"),
    @("", "19:21:32", "20:10", "23_HERNANDEZ",   49, "LP1912", "31/12/2025"),
    @("", "19:21:32", "20:14", "11_ETCHEVERRY",  53, "LP1912", "31/12/2025"),
    @("", "19:21:32", "20:26", "15_ABASTO",      65, "LP1912", "31/12/2025"),
    @("", "19:21:32", "20:28", "10_OLMOS",       67, "LP1912", "31/12/2025"),
    @("", "19:21:32", "20:44", "215B_EL PATO",   83, "LP1912", "31/12/2025"),
    @("", "19:21:32", "20:45", "17X38_ROMERO",   84, "LP1912", "31/12/2025"),
    @("", "19:21:32", "20:50", "23_HERNANDEZ",   89, "LP1912", "31/12/2025"),
    @("", "19:21:32", "20:53", "27_EL RETIRO",   92, "LP1912", "31/12/2025")
)

$r = 1258
foreach ($row in $rows1) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $ws1.Cells.Item($r, 6).Value = $row[5]
    $ws1.Cells.Item($r, 7).Value = $row[6]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215  (columns: A=meta, B=Fecha, C=Hora_Scrap, D=Hora_Llegada,
#                                 E=Linea, F=Minutos, G=Parada)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 31/12/2025 19:21:43"
$ws2.Range("A3").Value = "Total filas: 91"

$rows2 = @(
    @("", "31/12/2025", "19:21:32", "19:31", "215_EL PELIGRO", 10, "LP1912"),
    @("", "31/12/2025", "19:21:32", "20:01", "215C_EL PATO",   40, "LP1912"),
    @("", "31/12/2025", "19:21:32", "20:44", "215B_EL PATO",   83, "LP1912")
)

$r = 90
foreach ($row in $rows2) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
    $ws2.Cells.Item($r, 6).Value = $row[5]
    $ws2.Cells.Item($r, 7).Value = $row[6]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173  (columns: A=meta, B=Fecha, C=Hora_Scrap, D=Hora_Llegada,
#                                E=Linea, F=Minutos, G=Parada)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 31/12/2025 19:21:43"
$ws3.Range("A3").Value = "Total filas: 147"

$ws3.Cells.Item(148, 1).Value = ""
$ws3.Cells.Item(148, 2).Value = "31/12/2025"
$ws3.Cells.Item(148, 3).Value = "19:21:38"
$ws3.Cells.Item(148, 4).Value = "20:41"
$ws3.Cells.Item(148, 5).Value = "215C_LA PLATA"
$ws3.Cells.Item(148, 6).Value = 80
$ws3.Cells.Item(148, 7).Value = "L6203"
